$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# --- Row 12 (VT187-0208): fill in validate_Result with the captured image marker ---
$ws.Range("H12").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB and RE2.2 Semi Auto Frame Work : Imager`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0208`n};`nvalidate4`n{`nvalidate_Result=data : image/png;base64`n};"

# --- Row 13 (VT187-0211): fill in validate_Result with the network error code (fixes network issue) ---
$ws.Range("H13").Value = "validate1`n{`nvalidate_PageTitle=RE 2.2 Tests`n};`nvalidate2`n{`nvalidate_PageTitle=PB and RE2.2 Semi Auto Frame Work : Imager`n};`nvalidate3`n{`nvalidate_Text_Exists=VT187-0211`n};`nvalidate4`n{`nvalidate_Result=Error Code: 12014`n};"

# --- Row 15 (VT187-0224): extend the post-run wait from 5s to 15s ---
$ws.Range("G15").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nwait(5);`nSelectImager(back_camera);`nSelectTestToRun(VT187_0224_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"

# --- Row 16 (VT187-0235): extend the post-run wait from 5s to 15s ---
$ws.Range("G16").Value = "wait(3);`nvalidate1;`nlink_Click(imager_test_link);`nvalidate2;`nwait(5);`nSelectImager(back_camera);`nSelectTestToRun(VT187_0235_string);`nvalidate3;`nClickRunTest(runtest_top_xpath);`nwait(15);`nvalidate4;"

# --- Normalize formatting: D12, D13, D16 had a stray yellow-highlight style; match the rest of column D ---
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Record the final selection as it was when the workbook was last saved ---
$ws.Range("E2").Select()
